# Generate Report for Handoff
#
# The source file "3d63a2f5-6909-4861-82fd-5c387ddcbe82.md" finished localization
# handoff/handback (producing fresh .xlf artifacts and dropping the now-stale
# dependent .png rows), and a brand new source file
# "cccb50be-7698-471e-8f8f-50e3f32d2e44.md" was picked up and already handed off.
# This updates the Overview/zh-cn/de-de sheets accordingly.

$wb = $excel.ActiveWorkbook

$base_repo = "https://github.com/OpenLocalizationTest/oltest/blob/aea8dfa3f27b39cb6b9d008df46631b7c427e13c"
$zhcn_handoff_base = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9ae31ae3ea8b3bd93c0970c3f415031f80f1cd3e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$dede_handoff_base = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e7251a57269b746f5840bee5df16d8c138493c78/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$md1 = "acc5b70d-3d70-4ed8-81fb-134c95d14776.md"
$md2 = "cccb50be-7698-471e-8f8f-50e3f32d2e44.md"
$xlf1_zhcn = "acc5b70d-3d70-4ed8-81fb-134c95d14776.61a819e6173c7761cfbd7fb956b91698b7d75feb.zh-cn.xlf"
$xlf2_zhcn = "cccb50be-7698-471e-8f8f-50e3f32d2e44.1eb06eccd845a22c5bbe0bb21b016481da31cd8f.zh-cn.xlf"
$xlf1_dede = "acc5b70d-3d70-4ed8-81fb-134c95d14776.61a819e6173c7761cfbd7fb956b91698b7d75feb.de-de.xlf"
$xlf2_dede = "cccb50be-7698-471e-8f8f-50e3f32d2e44.1eb06eccd845a22c5bbe0bb21b016481da31cd8f.de-de.xlf"

$zhcn_time = "2016-03-08 23:26:30"
$dede_time = "2016-03-08 23:26:39"

# Removes every hyperlink on $sheet whose anchor Range.Address matches one of
# the addresses in $addrs. Re-scans the live collection from scratch before
# each single deletion because deleting mid-enumeration corrupts the
# Hyperlinks collection's enumerator.
function Remove-HyperlinksAt($sheet, $addrs) {
    foreach ($target in $addrs) {
        $keepGoing = $true
        while ($keepGoing) {
            $found = $null
            foreach ($hl in $sheet.Hyperlinks) {
                if ($hl.Range.Address() -eq $target) {
                    $found = $hl
                    break
                }
            }
            if ($found -eq $null) {
                $keepGoing = $false
            } else {
                $found.Delete()
            }
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.Address = "$base_repo/e2e/$md1"
        $hl.TextToDisplay = $md1
    } elseif ($addr -eq '$A$3') {
        $hl.Address = "$base_repo/e2e/$md2"
        $hl.TextToDisplay = $md2
    } elseif ($addr -eq '$A$4') {
        $hl.Address = "$base_repo/.localization-config"
        $hl.TextToDisplay = ".localization-config"
    }
}
Remove-HyperlinksAt $wsOverview @('$A$5')

$wsOverview.Range("A2").Value = $md1
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Range("A3").Value = $md2
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Range("A4").Value = ".localization-config"
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

$wsOverview.Range("A5:C5").Clear()

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.Address = "$base_repo/e2e/$md1"
        $hl.TextToDisplay = $md1
    } elseif ($addr -eq '$C$2') {
        $hl.Address = "$zhcn_handoff_base/$xlf1_zhcn"
        $hl.TextToDisplay = $xlf1_zhcn
    } elseif ($addr -eq '$A$3') {
        $hl.Address = "$base_repo/e2e/$md2"
        $hl.TextToDisplay = $md2
    } elseif ($addr -eq '$C$3') {
        $hl.Address = "$zhcn_handoff_base/$xlf2_zhcn"
        $hl.TextToDisplay = $xlf2_zhcn
    } elseif ($addr -eq '$A$4') {
        $hl.Address = "$base_repo/.localization-config"
        $hl.TextToDisplay = ".localization-config"
    }
}
Remove-HyperlinksAt $wsZhCn @('$C$4', '$A$5')

$wsZhCn.Range("A2").Value = $md1
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = $xlf1_zhcn
$wsZhCn.Range("D2").Value = $zhcn_time
$wsZhCn.Range("G2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H2").Value = "Include"

$wsZhCn.Range("A3").Value = $md2
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = $xlf2_zhcn
$wsZhCn.Range("D3").Value = $zhcn_time
$wsZhCn.Range("G3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H3").Value = "Include"
$wsZhCn.Range("I3").Clear()

$wsZhCn.Range("A4").Value = ".localization-config"
$wsZhCn.Range("B4").Value = "Not to be localized"
$wsZhCn.Range("C4").Clear()
$wsZhCn.Range("D4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("G4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H4").Value = "Ignored"
$wsZhCn.Range("I4").Clear()

$wsZhCn.Range("A5:I5").Clear()

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.Address = "$base_repo/e2e/$md1"
        $hl.TextToDisplay = $md1
    } elseif ($addr -eq '$C$2') {
        $hl.Address = "$dede_handoff_base/$xlf1_dede"
        $hl.TextToDisplay = $xlf1_dede
    } elseif ($addr -eq '$A$3') {
        $hl.Address = "$base_repo/e2e/$md2"
        $hl.TextToDisplay = $md2
    } elseif ($addr -eq '$C$3') {
        $hl.Address = "$dede_handoff_base/$xlf2_dede"
        $hl.TextToDisplay = $xlf2_dede
    } elseif ($addr -eq '$A$4') {
        $hl.Address = "$base_repo/.localization-config"
        $hl.TextToDisplay = ".localization-config"
    }
}
Remove-HyperlinksAt $wsDeDe @('$C$4', '$A$5')

$wsDeDe.Range("A2").Value = $md1
$wsDeDe.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = $xlf1_dede
$wsDeDe.Range("D2").Value = $dede_time
$wsDeDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H2").Value = "Include"

$wsDeDe.Range("A3").Value = $md2
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = $xlf2_dede
$wsDeDe.Range("D3").Value = $dede_time
$wsDeDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H3").Value = "Include"
$wsDeDe.Range("I3").Clear()

$wsDeDe.Range("A4").Value = ".localization-config"
$wsDeDe.Range("B4").Value = "Not to be localized"
$wsDeDe.Range("C4").Clear()
$wsDeDe.Range("D4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H4").Value = "Ignored"
$wsDeDe.Range("I4").Clear()

$wsDeDe.Range("A5:I5").Clear()
